$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "charge calculator" block (rows 125-130) ---
# Labels/units are entered first, in the order that reproduces the
# shared-strings table growth seen in the target workbook.
$ws.Range("I125").Value = "q = I * t"
$ws.Range("E126").Value = "t="
$ws.Range("E130").Value = "I="
$ws.Range("G130").Value = "A"
$ws.Range("H126").Value = "q"
$ws.Range("J126").Value = "C"
$ws.Range("J127").Value = "mC"
$ws.Range("J128").Value = "uC"
$ws.Range("J129").Value = "nC"

# Row 125: frequency-like input value and its label
$ws.Range("F125").Value = 300000

# Row 126: t = 1/F125 (seconds), unit "s", q = F126*F130 (Coulombs)
$ws.Range("F126").Formula = "=1/F125"
$ws.Range("G126").Value = "s"
$ws.Range("I126").Formula = "=F126*F130"

# Row 127: t in ms, q in mC
$ws.Range("F127").Formula = "=F126*1000"
$ws.Range("G127").Value = "ms"
$ws.Range("I127").Formula = "=I126*1000"

# Row 128: t in us, q in uC
$ws.Range("F128").Formula = "=F127*1000"
$ws.Range("G128").Value = "us"
$ws.Range("I128").Formula = "=I127*1000"

# Row 129: q in nC
$ws.Range("I129").Formula = "=I128*1000"

# Row 130: I = 0.018 A
$ws.Range("F130").Value = 0.018

# Restore the selection to where the author left off editing
$ws.Range("F126").Select()
